# Auto-generated edit script: applies numeric cell updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# matching the scheduled-runner data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 60420
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 60420
$ws.Range("K70").Value = 0
$ws.Range("N70").Value = -181800
$ws.Range("L70").Value = 181260
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 60420
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 60420
$ws.Range("K73").Value = 0
$ws.Range("N73").Value = -183132
$ws.Range("L73").Value = 181260
$ws.Range("M73").ClearContents()

$ws.Range("H116").Value = 21220948
$ws.Range("I116").Value = 16668766
$ws.Range("J116").Value = 23822196
$ws.Range("K116").Value = 16668766
$ws.Range("L116").Value = 23822196
$ws.Range("M116").Value = -16665324
$ws.Range("N116").Value = -23829080

$ws.Range("H132").Value = 4116538.2
$ws.Range("I132").Value = 1362.4615
$ws.Range("K132").Value = 4087.3845
$ws.Range("M132").Value = -1557.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2872145.2
$ws.Range("I61").Value = 1389927.6
$ws.Range("J61").Value = 11765451
$ws.Range("K61").Value = 1389927.6
$ws.Range("L61").Value = 11765451
$ws.Range("M61").Value = -1389715.6
$ws.Range("N61").Value = -11765875

$ws.Range("H74").Value = 81454720
$ws.Range("I74").Value = 67619580
$ws.Range("J74").Value = 133336500
$ws.Range("K74").Value = 67619580
$ws.Range("L74").Value = 133336500
$ws.Range("M74").Value = -67618706
$ws.Range("N74").Value = -133338248

$ws.Range("H77").Value = 81454720
$ws.Range("I77").Value = 67619580
$ws.Range("J77").Value = 133336500
$ws.Range("K77").Value = 338097900
$ws.Range("L77").Value = 666682500
$ws.Range("M77").Value = -338093532
$ws.Range("N77").Value = -666691236

$ws.Range("H132").Value = 16209097
$ws.Range("I132").Value = 15878703
$ws.Range("K132").Value = 47636109
$ws.Range("M132").Value = -47633579

$ws.Range("H136").Value = 2872145.2
$ws.Range("I136").Value = 1389927.6
$ws.Range("J136").Value = 11765451
$ws.Range("K136").Value = 4169782.8
$ws.Range("L136").Value = 35296353
$ws.Range("M136").Value = -4167232.8
$ws.Range("N136").Value = -35301453

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8504736
$ws.Range("I134").Value = 9260354
$ws.Range("J134").Value = 3971029.5
$ws.Range("K134").Value = 27781062
$ws.Range("L134").Value = 11913088.5
$ws.Range("M134").Value = -27778527
$ws.Range("N134").Value = -11918158.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3160518.2
$ws.Range("J31").Value = 3682098.8
$ws.Range("L31").Value = 3682098.8
$ws.Range("N31").Value = -3682688.8

$ws.Range("H34").Value = 3160518.2
$ws.Range("J34").Value = 3682098.8
$ws.Range("L34").Value = 3682098.8
$ws.Range("N34").Value = -3682502.8

$ws.Range("H58").Value = 2389936.2
$ws.Range("I58").Value = 1556090.4
$ws.Range("K58").Value = 1556090.4
$ws.Range("M58").Value = -1555887.4

$ws.Range("H62").Value = 2700
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 2700
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 2700
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -3948

$ws.Range("H65").Value = 2700
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 2700
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 13500
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -19740

$ws.Range("H86").Value = 8337.857
$ws.Range("I86").Value = 13223.454
$ws.Range("J86").Value = 5176.5884
$ws.Range("K86").Value = 13223.454
$ws.Range("L86").Value = 5176.5884
$ws.Range("M86").Value = -12100.454
$ws.Range("N86").Value = -7422.5884

$ws.Range("H89").Value = 8337.857
$ws.Range("I89").Value = 13223.454
$ws.Range("J89").Value = 5176.5884
$ws.Range("K89").Value = 66117.27
$ws.Range("L89").Value = 25882.942
$ws.Range("M89").Value = -60501.27
$ws.Range("N89").Value = -37114.942

$ws.Range("H99").Value = 12927.889
$ws.Range("I99").Value = 10787.5
$ws.Range("J99").Value = 14640.2
$ws.Range("K99").Value = 10787.5
$ws.Range("L99").Value = 14640.2
$ws.Range("M99").Value = -9289.5
$ws.Range("N99").Value = -17636.2

$ws.Range("H105").Value = 8043.5
$ws.Range("I105").Value = 1719.6
$ws.Range("K105").Value = 1719.6
$ws.Range("M105").Value = 27.40000000000009

$ws.Range("H107").Value = 687.6875
$ws.Range("I107").Value = 315
$ws.Range("J107").Value = 911.3
$ws.Range("K107").Value = 315
$ws.Range("L107").Value = 911.3
$ws.Range("M107").Value = 1605
$ws.Range("N107").Value = -4751.3

$ws.Range("H122").Value = 4794.7393
$ws.Range("I122").Value = 9291.091
$ws.Range("J122").Value = 673.0833
$ws.Range("K122").Value = 27873.273
$ws.Range("L122").Value = 2019.2499
$ws.Range("M122").Value = -25423.273
$ws.Range("N122").Value = -6919.2499

$ws.Range("H126").Value = 12927.889
$ws.Range("I126").Value = 10787.5
$ws.Range("J126").Value = 14640.2
$ws.Range("K126").Value = 32362.5
$ws.Range("L126").Value = 43920.60000000001
$ws.Range("M126").Value = -29892.5
$ws.Range("N126").Value = -48860.60000000001

$ws.Range("H132").Value = 2275656.2
$ws.Range("I132").Value = 3335361.5
$ws.Range("K132").Value = 10006084.5
$ws.Range("M132").Value = -10003554.5

$ws.Range("H134").Value = 1295849.6
$ws.Range("I134").Value = 5454
$ws.Range("J134").Value = 5720063
$ws.Range("K134").Value = 16362
$ws.Range("L134").Value = 17160189
$ws.Range("M134").Value = -13827
$ws.Range("N134").Value = -17165259

$ws.Range("H136").Value = 2389936.2
$ws.Range("I136").Value = 1556090.4
$ws.Range("K136").Value = 4668271.199999999
$ws.Range("M136").Value = -4665721.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2436.5557
$ws.Range("I9").Value = 1500
$ws.Range("J9").Value = 2491.647
$ws.Range("K9").Value = 4500
$ws.Range("L9").Value = 7474.941
$ws.Range("M9").Value = -4276
$ws.Range("N9").Value = -7922.941

$ws.Range("H10").Value = 257.5
$ws.Range("I10").Value = 99.875
$ws.Range("J10").Value = 888
$ws.Range("K10").Value = 299.625
$ws.Range("L10").Value = 2664
$ws.Range("M10").Value = -160.625
$ws.Range("N10").Value = -2942

$ws.Range("H113").Value = 1141.909
$ws.Range("J113").Value = 2034.8695
$ws.Range("L113").Value = 6104.6085
$ws.Range("N113").Value = -10444.6085

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2488376
$ws.Range("I70").Value = 1592981.1
$ws.Range("J70").Value = 3794160.5
$ws.Range("K70").Value = 1592981.1
$ws.Range("L70").Value = 3794160.5
$ws.Range("M70").Value = -1592711.1
$ws.Range("N70").Value = -3794700.5

$ws.Range("H73").Value = 2488376
$ws.Range("I73").Value = 1592981.1
$ws.Range("J73").Value = 3794160.5
$ws.Range("K73").Value = 1592981.1
$ws.Range("L73").Value = 3794160.5
$ws.Range("M73").Value = -1592045.1
$ws.Range("N73").Value = -3796032.5

$ws.Range("H97").Value = 13159565
$ws.Range("I97").Value = 1009.9286
$ws.Range("J97").Value = 50003520
$ws.Range("K97").Value = 1009.9286
$ws.Range("L97").Value = 50003520
$ws.Range("M97").Value = -513.9286
$ws.Range("N97").Value = -50004512

$ws.Range("H113").Value = 26808
$ws.Range("I113").Value = 7129.9
$ws.Range("J113").Value = 76003.25
$ws.Range("K113").Value = 7129.9
$ws.Range("L113").Value = 76003.25
$ws.Range("M113").Value = -4959.9
$ws.Range("N113").Value = -80343.25

$ws.Range("H132").Value = 56424356
$ws.Range("I132").Value = 123809520
$ws.Range("J132").Value = 22731772
$ws.Range("K132").Value = 371428560
$ws.Range("L132").Value = 68195316
$ws.Range("M132").Value = -371426030
$ws.Range("N132").Value = -68200376

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1668608
$ws.Range("I132").Value = 2470570.8
$ws.Range("J132").Value = 2993.1538
$ws.Range("K132").Value = 7411712.399999999
$ws.Range("L132").Value = 8979.4614
$ws.Range("M132").Value = -7409182.399999999
$ws.Range("N132").Value = -14039.4614

$ws.Range("H136").Value = 1471857.8
$ws.Range("I136").Value = 2102031.2
$ws.Range("J136").Value = 1452.5834
$ws.Range("K136").Value = 6306093.600000001
$ws.Range("L136").Value = 4357.7502
$ws.Range("M136").Value = -6303543.600000001
$ws.Range("N136").Value = -9457.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1410872.4
$ws.Range("I132").Value = 1016529.94
$ws.Range("J132").Value = 2528176
$ws.Range("K132").Value = 3049589.82
$ws.Range("L132").Value = 7584528
$ws.Range("M132").Value = -3047059.82
$ws.Range("N132").Value = -7589588

$ws.Range("H136").Value = 12408.55
$ws.Range("I136").Value = 8457.4
$ws.Range("K136").Value = 25372.2
$ws.Range("M136").Value = -22822.2
